# Daily attendance processing - 2026-02-21 20:33:30 UTC
# Reorders the comma-separated "Recorded By" values (column G) for specific
# rows so the later-added/most-recent recorder (2025/2026 session or the
# session owner) appears first in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    3  = "2022/2023, 2025/2026"
    22 = "2025/2026, 2024/2025"
    23 = "2022/2023, 2025/2026, 2023/2024"
    24 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    27 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    28 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    31 = "2022/2023, 2025/2026"
    50 = "2025/2026, 2024/2025"
    51 = "2022/2023, 2025/2026, 2023/2024"
    52 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    55 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    56 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
